# Apply scraped market-data refresh to Sheets (rows updated by scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 5628.0386
$ws.Range("I86").Value = 15359.286
$ws.Range("J86").Value = 2042.8422
$ws.Range("K86").Value = 15359.286
$ws.Range("L86").Value = 2042.8422
$ws.Range("M86").Value = -14236.286
$ws.Range("N86").Value = -4288.8422

# Row 89
$ws.Range("H89").Value = 5628.0386
$ws.Range("I89").Value = 15359.286
$ws.Range("J89").Value = 2042.8422
$ws.Range("K89").Value = 76796.42999999999
$ws.Range("L89").Value = 10214.211
$ws.Range("M89").Value = -71180.42999999999
$ws.Range("N89").Value = -21446.211

# Row 99
$ws.Range("H99").Value = 685.875
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 9000
$ws.Range("N99").Value = -11996

# Row 138
$ws.Range("H138").Value = 3475.638
$ws.Range("I138").Value = 1638.1034
$ws.Range("J138").Value = 5313.1724
$ws.Range("K138").Value = 4914.3102
$ws.Range("L138").Value = 15939.5172
$ws.Range("M138").Value = 225.6898000000001
$ws.Range("N138").Value = -26219.5172

# Row 141
$ws.Range("H141").Value = 619770.4
$ws.Range("I141").Value = 2306
$ws.Range("J141").Value = 1391600.9
$ws.Range("K141").Value = 6918
$ws.Range("L141").Value = 4174802.7
$ws.Range("M141").Value = -1738
$ws.Range("N141").Value = -4185162.7

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 505000
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10288

# Row 61
$ws.Range("H61").Value = 5750
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1288
$ws.Range("N61").Value = -10424

# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

# Row 75
$ws.Range("H75").Value = 29000
$ws.Range("J75").Value = 29000
$ws.Range("L75").Value = 29000
$ws.Range("N75").Value = -30748

# Row 78
$ws.Range("H78").Value = 29000
$ws.Range("J78").Value = 29000
$ws.Range("L78").Value = 87000
$ws.Range("N78").Value = -95736

# Row 81
$ws.Range("H81").Value = 29590.5
$ws.Range("J81").Value = 29590.5
$ws.Range("L81").Value = 29590.5
$ws.Range("N81").Value = -31586.5

# Row 82
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 84
$ws.Range("H84").Value = 29590.5
$ws.Range("J84").Value = 29590.5
$ws.Range("L84").Value = 88771.5
$ws.Range("N84").Value = -98755.5

# Row 85
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 86
$ws.Range("H86").Value = 40000
$ws.Range("J86").Value = 40000
$ws.Range("L86").Value = 40000
$ws.Range("N86").Value = -42372

# Row 89
$ws.Range("H89").Value = 40000
$ws.Range("J89").Value = 40000
$ws.Range("L89").Value = 120000
$ws.Range("N89").Value = -131856

# Row 107
$ws.Range("H107").Value = 23000
$ws.Range("J107").Value = 23000
$ws.Range("L107").Value = 23000
$ws.Range("N107").Value = -30680

# Row 136
$ws.Range("H136").Value = 5750
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2469.1428
$ws.Range("I99").Value = 1342.5454
$ws.Range("K99").Value = 1342.5454
$ws.Range("M99").Value = 155.4546

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 59500.75
$ws.Range("I3").Value = 30000
$ws.Range("J3").Value = 69334.336
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 69334.336
$ws.Range("M3").Value = -29887
$ws.Range("N3").Value = -69560.336

# Row 132
$ws.Range("H132").Value = 3148.9697
$ws.Range("I132").Value = 1714.45
$ws.Range("J132").Value = 5355.923
$ws.Range("K132").Value = 5143.35
$ws.Range("L132").Value = 16067.769
$ws.Range("M132").Value = -2613.35
$ws.Range("N132").Value = -21127.769

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

# Row 138
$ws.Range("H138").Value = 29211.3
$ws.Range("J138").Value = 29211.3
$ws.Range("L138").Value = 29211.3
$ws.Range("N138").Value = -39491.3

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 877.89655
$ws.Range("I131").Value = 440.41666
$ws.Range("J131").Value = 2977.8
$ws.Range("K131").Value = 1321.24998
$ws.Range("L131").Value = 8933.400000000001
$ws.Range("M131").Value = 3718.75002
$ws.Range("N131").Value = -19013.4

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 20224.25
$ws.Range("I3").Value = 297.66666
$ws.Range("J3").Value = 80004
$ws.Range("K3").Value = 297.66666
$ws.Range("L3").Value = 80004
$ws.Range("M3").Value = -181.66666
$ws.Range("N3").Value = -80236

# Row 43
$ws.Range("H43").Value = 4499
$ws.Range("I43").Value = 3495
$ws.Range("J43").Value = 4750
$ws.Range("K43").Value = 3495
$ws.Range("L43").Value = 4750
$ws.Range("M43").Value = -3344
$ws.Range("N43").Value = -5052

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2199
$ws.Range("I7").Value = 1369.5714
$ws.Range("J7").Value = 3166.6667
$ws.Range("K7").Value = 1369.5714
$ws.Range("L7").Value = 3166.6667
$ws.Range("M7").Value = -1257.5714
$ws.Range("N7").Value = -3390.6667

# Row 16
$ws.Range("H16").Value = 8983.333000000001
$ws.Range("I16").Value = 2900
$ws.Range("J16").Value = 10200
$ws.Range("K16").Value = 2900
$ws.Range("L16").Value = 10200
$ws.Range("M16").Value = -2730
$ws.Range("N16").Value = -10540

# Row 40
$ws.Range("H40").Value = 2499.6667
$ws.Range("J40").Value = 4399.3335
$ws.Range("L40").Value = 4399.3335
$ws.Range("N40").Value = -4671.3335

# Row 126
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 1369.5714
$ws.Range("J126").Value = 3166.6667
$ws.Range("K126").Value = 4108.7142
$ws.Range("L126").Value = 9500.000100000001
$ws.Range("M126").Value = -1638.7142
$ws.Range("N126").Value = -14440.0001

# Row 132
$ws.Range("H132").Value = 2608.2083
$ws.Range("I132").Value = 2020.4
$ws.Range("J132").Value = 3028.0715
$ws.Range("K132").Value = 6061.200000000001
$ws.Range("L132").Value = 9084.2145
$ws.Range("M132").Value = -3531.200000000001
$ws.Range("N132").Value = -14144.2145

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 318013.22
$ws.Range("I132").Value = 527990.25
$ws.Range("J132").Value = 11123.692
$ws.Range("K132").Value = 1583970.75
$ws.Range("L132").Value = 33371.076
$ws.Range("M132").Value = -1581440.75
$ws.Range("N132").Value = -38431.076

